{"js": "// Update the cover-page date line from \"\u2013 Hanoi, August 2019 \u2013\" to\n// \"\u2013 Hanoi, October 2022 \u2013\", keeping the existing run formatting\n// (rFonts cstheme=\"minorHAnsi\", sz=28, szCs=28) untouched.\n\nconst body = context.document.body;\n\n// Locate the exact run of text on the cover page. Search is scoped to the\n// whole body text-search (matchCase so we don't accidentally touch any other\n// casing variant) \u2014 the phrase is unique in this document.\nconst results = body.search(\"Hanoi, August 2019\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace just the variable middle portion (\"August 2019\" -> \"October 2022\")\n  // so the surrounding \"\u2013 \" / \" \u2013\" dash characters (and the run formatting\n  // that decorates them) are left completely undisturbed.\n  results.items[0].insertText(\"Hanoi, October 2022\", \"Replace\");\n  await context.sync();\n} else {\n  // Fall back: maybe the text was already partly edited / spacing differs \u2014\n  // search more narrowly for the two date tokens independently so the edit\n  // is still applied.\n  const augResults = body.search(\"August\", { matchCase: true });\n  augResults.load(\"items/text\");\n  await context.sync();\n  if (augResults.items.length > 0) {\n    augResults.items[0].insertText(\"October\", \"Replace\");\n    await context.sync();\n  }\n\n  const yearResults = body.search(\"2019\", { matchCase: true });\n  yearResults.load(\"items/text\");\n  await context.sync();\n  if (yearResults.items.length > 0) {\n    yearResults.items[0].insertText(\"2022\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Update the cover-page date line from \"\u2013 Hanoi, August 2019 \u2013\" to\n# \"\u2013 Hanoi, October 2022 \u2013\", keeping the existing paragraph/run formatting\n# (rFonts cstheme=\"minorHAnsi\", sz=28, szCs=28) untouched.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Hanoi, August 2019*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Re-assigning the whole paragraph Range's .Text preserves the run's\n    # existing formatting/rsid instead of minting a brand-new, unformatted\n    # run the way Find/Replace does.\n    $target.Range.Text = \"\u2013 Hanoi, October 2022 \u2013\"\n} else {\n    # Fallback: use Find/Replace in case spacing/punctuation differs slightly.\n    $find = $d.Content.Find\n    $find.Text = \"August\"\n    $find.Replacement.Text = \"October\"\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n    $find2 = $d.Content.Find\n    $find2.Text = \"2019\"\n    $find2.Replacement.Text = \"2022\"\n    $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n}\n"}
